$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Row 24: D30 entry
$ws.Range("B24").Value = "D30"
$ws.Range("C24").Value = "2/3/2020"
$ws.Range("C24").NumberFormat = "mm-dd-yy"
$ws.Range("D24").Value = "Started with lesson 8 - reached 8.3"

# Row 25: D31 entry
$ws.Range("B25").Value = "D31"
$ws.Range("C25").Value = "2/4/2020"
$ws.Range("C25").NumberFormat = "mm-dd-yy"
$ws.Range("D25").Value = "Day off"

# Update the view state to match the diff (scrolled position & selection)
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("E26").Select()
